# Auto-generated Excel COM-interop script
# Applies the diff: numeric corrections across sheets + 4 new rows on 展览

$wb = $excel.ActiveWorkbook

# ---- 展览: numeric corrections (rows 1-40) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value = 2819
$ws1.Cells.Item(5, 6).Value = 6380
$ws1.Cells.Item(6, 6).Value = 2465
$ws1.Cells.Item(8, 6).Value = 32
$ws1.Cells.Item(9, 6).Value = 36
$ws1.Cells.Item(10, 6).Value = 2886
$ws1.Cells.Item(12, 6).Value = 29
$ws1.Cells.Item(13, 6).Value = 7056
$ws1.Cells.Item(14, 6).Value = 283
$ws1.Cells.Item(16, 6).Value = 218
$ws1.Cells.Item(19, 6).Value = 8217
$ws1.Cells.Item(21, 6).Value = 39
$ws1.Cells.Item(27, 6).Value = 67
$ws1.Cells.Item(33, 6).Value = 2596
$ws1.Cells.Item(34, 6).Value = 41
$ws1.Cells.Item(37, 6).Value = 1156
$ws1.Cells.Item(39, 6).Value = 654
$ws1.Cells.Item(40, 6).Value = 3660
$ws1.Cells.Item(3, 7).Value = 39.9

# ---- 展览: rewrite rows 41-43 + append new rows 44-47 ----
# Helper cell (far outside the table) used to force literal-text assignment
# for date-like strings that Excel would otherwise auto-convert to dates.
$helper = $ws1.Range("A500")
function Set-TextValue($range, [string]$text) {
    $helper.NumberFormat = "@"
    $helper.Value = $text
    $helper.Copy()
    $range.PasteSpecial(-4163)
}

# Give the new rows 44-47 the same formatting (bold/border/center) as column A
# already has on the existing data rows.
$ws1.Range("A43").Copy()
$ws1.Range("A44:A47").PasteSpecial(-4122)

# Row 41: 北京·第十三届GOJO超次元动漫游戏嘉年华·一周年盛典
$ws1.Range("A41").Value = 40
Set-TextValue $ws1.Range("B41") "2024-06-01"
$ws1.Range("C41").Value = "北京·第十三届GOJO超次元动漫游戏嘉年华·一周年盛典"
$ws1.Range("D41").Value = "小关路39号 北投购物公园"
Set-TextValue $ws1.Range("E41") "2024.06.01 10:00-06.02 17:00"
$ws1.Range("F41").Value = 1
$ws1.Range("G41").Value = 6.6
$ws1.Range("H41").Value = "https://show.bilibili.com/platform/detail.html?id=83827"
$ws1.Range("I41").Value = "//i2.hdslb.com/bfs/openplatform/202404/A4pgKCpJ1712043154165.jpeg"

# Row 42: 北京·ICOS国际动漫节×CGF中国游戏节02
$ws1.Range("A42").Value = 41
Set-TextValue $ws1.Range("B42") "2024-06-08"
$ws1.Range("C42").Value = "北京·ICOS国际动漫节×CGF中国游戏节02"
$ws1.Range("D42").Value = "石景山路68号 北京首钢会展中心"
Set-TextValue $ws1.Range("E42") "2024.06.08 09:00-06.09 17:00"
$ws1.Range("F42").Value = 175
$ws1.Range("G42").Value = 70
$ws1.Range("H42").Value = "https://show.bilibili.com/platform/detail.html?id=83161"
$ws1.Range("I42").Value = "//i2.hdslb.com/bfs/openplatform/202403/965YPK7G1711003390117.jpeg"

# Row 43: 北京·thebONE游戏动漫节
$ws1.Range("A43").Value = 42
Set-TextValue $ws1.Range("B43") "2024-06-08"
$ws1.Range("C43").Value = "北京·thebONE游戏动漫节"
$ws1.Range("D43").Value = "小关路39号 北投购物公园"
Set-TextValue $ws1.Range("E43") "2024.06.08 10:00-06.10 17:00"
$ws1.Range("F43").Value = 3
$ws1.Range("G43").Value = 75
$ws1.Range("H43").Value = "https://show.bilibili.com/platform/detail.html?id=83830"
$ws1.Range("I43").Value = "//i0.hdslb.com/bfs/openplatform/202404/PAQ2DFrV1712046388743.jpeg"

# Row 44: 北京·万游引力国潮动漫嘉年华s7
$ws1.Range("A44").Value = 43
Set-TextValue $ws1.Range("B44") "2024-06-22"
$ws1.Range("C44").Value = "北京·万游引力国潮动漫嘉年华s7"
$ws1.Range("D44").Value = "北七家镇王府街55号 水城会议中心"
Set-TextValue $ws1.Range("E44") "2024.06.22 10:00-06.23 17:00"
$ws1.Range("F44").Value = 1177
$ws1.Range("G44").Value = 75
$ws1.Range("H44").Value = "https://show.bilibili.com/platform/detail.html?id=82848"
$ws1.Range("I44").Value = "//i2.hdslb.com/bfs/openplatform/202403/DFRkPH7q1710396818728.jpeg"

# Row 45: 北京·IDO动漫游戏嘉年华46th
$ws1.Range("A45").Value = 44
Set-TextValue $ws1.Range("B45") "2024-07-20"
$ws1.Range("C45").Value = "北京·IDO动漫游戏嘉年华46th"
$ws1.Range("D45").Value = "京沈路与天北路交汇处西北角 中国国际展览中心新馆"
Set-TextValue $ws1.Range("E45") "2024.07.20 09:30-07.21 17:00"
$ws1.Range("F45").Value = 156
$ws1.Range("G45").Value = 75
$ws1.Range("H45").Value = "https://show.bilibili.com/platform/detail.html?id=83716"
$ws1.Range("I45").Value = "//i0.hdslb.com/bfs/openplatform/202404/G4DiYbc51712040520493.jpeg"

# Row 46: 北京·梦次元动漫展
$ws1.Range("A46").Value = 45
Set-TextValue $ws1.Range("B46") "2024-08-10"
$ws1.Range("C46").Value = "北京·梦次元动漫展"
$ws1.Range("D46").Value = "北京展览馆 北京展览馆"
Set-TextValue $ws1.Range("E46") "2024.08.10 10:00-08.11 17:00"
$ws1.Range("F46").Value = 3
$ws1.Range("G46").Value = 70
$ws1.Range("H46").Value = "https://show.bilibili.com/platform/detail.html?id=83828"
$ws1.Range("I46").Value = "//i1.hdslb.com/bfs/openplatform/202404/aUYvg6Cu1712054086278.jpeg"

# Row 47: 北京·IDO动漫游戏嘉年华47th
$ws1.Range("A47").Value = 46
Set-TextValue $ws1.Range("B47") "2024-10-01"
$ws1.Range("C47").Value = "北京·IDO动漫游戏嘉年华47th"
$ws1.Range("D47").Value = "亦庄荣昌东街6号 北京亦创国际会展中心"
Set-TextValue $ws1.Range("E47") "2024.10.01 10:00-10.03 17:00"
$ws1.Range("F47").Value = 3
$ws1.Range("G47").Value = 75
$ws1.Range("H47").Value = "https://show.bilibili.com/platform/detail.html?id=83826"
$ws1.Range("I47").Value = "//i1.hdslb.com/bfs/openplatform/202404/wxWMaLKJ1712054345299.jpeg"

$helper.Clear()

# ---- 演出: numeric correction ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(5, 6).Value = 244
$ws2.Cells.Item(5, 7).Value = 180

# ---- 全部类型: numeric corrections ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3, 6).Value = 2819
$ws4.Cells.Item(5, 6).Value = 244
$ws4.Cells.Item(6, 6).Value = 244
$ws4.Cells.Item(7, 6).Value = 6380
$ws4.Cells.Item(8, 6).Value = 2465
$ws4.Cells.Item(11, 6).Value = 32
$ws4.Cells.Item(12, 6).Value = 36
$ws4.Cells.Item(13, 6).Value = 2886
$ws4.Cells.Item(17, 6).Value = 29
$ws4.Cells.Item(18, 6).Value = 7056
$ws4.Cells.Item(19, 6).Value = 283
$ws4.Cells.Item(21, 6).Value = 218
$ws4.Cells.Item(24, 6).Value = 8217
$ws4.Cells.Item(26, 6).Value = 39
$ws4.Cells.Item(32, 6).Value = 67
$ws4.Cells.Item(39, 6).Value = 2596
$ws4.Cells.Item(40, 6).Value = 41
$ws4.Cells.Item(43, 6).Value = 1156
$ws4.Cells.Item(44, 6).Value = 654
$ws4.Cells.Item(46, 6).Value = 3660
$ws4.Cells.Item(47, 6).Value = 175
$ws4.Cells.Item(49, 6).Value = 1177
$ws4.Cells.Item(50, 6).Value = 156
$ws4.Cells.Item(3, 7).Value = 39.9
$ws4.Cells.Item(5, 7).Value = 180
$ws4.Cells.Item(6, 7).Value = 180

